# Update the lattice-multiplication exercise table: every cell keeps its
# existing "problem card" layout (operands / multiplier digits / dashed
# rule / lattice rows) but gets new numbers, per the target diff.
#
# Row/Col -> new 5-line cell content (joined later with a vertical-tab,
# which Word's Range.Text setter turns back into <w:t>/<w:br/> runs).
$cellValues = @(
    @(1, 1, @("63 x 73", "  7    3", "  ----", "6|    |", "3|    |")),
    @(1, 2, @("59 x 45", "  4    5", "  ----", "5|    |", "9|    |")),
    @(1, 3, @("75 x 81", "  8    1", "  ----", "7|    |", "5|    |")),
    @(2, 1, @("22 x 25", "  2    5", "  ----", "2|    |", "2|    |")),
    @(2, 2, @("54 x 23", "  2    3", "  ----", "5|    |", "4|    |")),
    @(2, 3, @("43 x 84", "  8    4", "  ----", "4|    |", "3|    |")),
    @(3, 1, @("85 x 53", "  5    3", "  ----", "8|    |", "5|    |")),
    @(3, 2, @("13 x 83", "  8    3", "  ----", "1|    |", "3|    |")),
    @(3, 3, @("71 x 49", "  4    9", "  ----", "7|    |", "1|    |")),
    @(4, 1, @("45 x 33", "  3    3", "  ----", "4|    |", "5|    |")),
    @(4, 2, @("15 x 73", "  7    3", "  ----", "1|    |", "5|    |")),
    @(4, 3, @("28 x 28", "  2    8", "  ----", "2|    |", "8|    |")),
    @(5, 1, @("98 x 88", "  8    8", "  ----", "9|    |", "8|    |")),
    @(5, 2, @("32 x 43", "  4    3", "  ----", "3|    |", "2|    |")),
    @(5, 3, @("40 x 20", "  2    0", "  ----", "4|    |", "0|    |"))
)

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$vt = [char]11

foreach ($entry in $cellValues) {
    $rowIdx = $entry[0]
    $colIdx = $entry[1]
    $linesArr = $entry[2]
    $newText = [string]::Join($vt, $linesArr)

    $cell = $t.Cell($rowIdx, $colIdx)
    $cell.Range.Text = $newText
}
